$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump the published version and timestamp.
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# Insert a new "Jurisdiction" metadata row right before "Description" (old row 11).
$ws.Rows.Item(11).Insert()

# Inherit the body-row formatting (border/alignment/style) from the row
# that just got pushed down, rather than whatever Insert() guessed.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
